$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update prior variance column (E2:E22) from 0.05 to 0.15
$ws.Range("E2:E22").Value = 0.15

# Move the active selection to E17 (also drops the stale topLeftCell scroll state)
$ws.Range("E17").Select()

# The workbook was switched to manual calculation for this round of test runs
$excel.Calculation = -4135  # xlCalculationManual
